# Applies the "Updated symbol list" data refresh (Tue Dec 27 09:07:43 UTC 2022)
# to the cryptos worksheet: for every data row (2-51) the "Hora" (hour) column
# moves from 8 to 9 and most "Price" values are refreshed; for rows 18-24 and
# 42-43 the coin list also shifted by one position, so Coin/Link/Price/Volume
# columns are updated to reflect the new row contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @{ row=2; D="242.53"; G="9" },
    @{ row=3; D="23.00"; G="9" },
    @{ row=4; D="5.410"; G="9" },
    @{ row=5; D="0.05954"; G="9" },
    @{ row=6; D="3.429"; G="9" },
    @{ row=7; D="6.499"; G="9" },
    @{ row=8; D="0.8145"; G="9" },
    @{ row=9; D="0.9242"; G="9" },
    @{ row=10; D="0.1437"; G="9" },
    @{ row=11; D="0.07413"; G="9" },
    @{ row=12; D="0.03295"; G="9" },
    @{ row=13; D="0.03076"; G="9" },
    @{ row=14; D="0.09354"; G="9" },
    @{ row=15; D="3.858"; G="9" },
    @{ row=16; D="0.001573"; G="9" },
    @{ row=17; D="0.04697"; G="9" },
    @{ row=18; B="One"; C="https://coinranking.com/coin/6Lga5NiXX3rT+one-one"; D="0.0005902"; E="17OneONE"; G="9" },
    @{ row=19; B="TigerCash"; C="https://coinranking.com/coin/6hIn06L2+tigercash-tch"; D="0.005864"; E="18TigerCashTCH"; G="9" },
    @{ row=20; B="BitKan"; C="https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"; D="0.001264"; E="19BitKanKAN"; G="9" },
    @{ row=21; B="HotbitToken"; C="https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"; D="0.004838"; E="20HotbitTokenHTB"; G="9" },
    @{ row=22; B="NitroEx"; C="https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"; D="0.00008004"; E="21NitroExNTXWorstin24h"; G="9" },
    @{ row=23; B="LEO"; C="https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; D="3.574"; E="22LEOLEO"; G="9" },
    @{ row=24; B="BTSEToken"; C="https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"; D="2.158"; E="23BTSETokenBTSE"; G="9" },
    @{ row=25; D="0.3238"; G="9" },
    @{ row=26; D="0.1331"; G="9" },
    @{ row=27; D="0.0002340"; G="9" },
    @{ row=28; G="9" },
    @{ row=29; G="9" },
    @{ row=30; G="9" },
    @{ row=31; G="9" },
    @{ row=32; G="9" },
    @{ row=33; G="9" },
    @{ row=34; G="9" },
    @{ row=35; G="9" },
    @{ row=36; G="9" },
    @{ row=37; G="9" },
    @{ row=38; G="9" },
    @{ row=39; G="9" },
    @{ row=40; D="0.03942"; G="9" },
    @{ row=41; D="0.006403"; G="9" },
    @{ row=42; B="BKEXToken"; C="https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"; D="0.1074"; E="41BKEXTokenBKK"; G="9" },
    @{ row=43; B="CEJI"; C="https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"; D="0.002601"; E="42CEJICEJI"; G="9" },
    @{ row=44; D="0.008909"; G="9" },
    @{ row=45; D="0.00005178"; G="9" },
    @{ row=46; G="9" },
    @{ row=47; D="0.7002"; G="9" },
    @{ row=48; D="0.002144"; G="9" },
    @{ row=49; D="0.00002101"; G="9" },
    @{ row=50; D="0.0002001"; G="9" },
    @{ row=51; G="9" }
)

# Columns that hold numeric-looking text (Price, Hora) must stay stored as
# text so values such as "23.00" or "0.0002340" keep their exact formatting
# instead of being coerced into numbers. Force the "Text" number format
# before writing, then drop back to the default "Normal" style so no new
# cell formatting is introduced.
$textColumns = @("D", "G")

foreach ($item in $rowData) {
    $r = $item.row
    foreach ($col in @("B", "C", "D", "E", "G")) {
        if ($item.ContainsKey($col)) {
            $cell = $ws.Range("$col$r")
            if ($textColumns -contains $col) {
                $cell.NumberFormat = "@"
                $cell.Value = $item[$col]
                $cell.Style = "Normal"
            } else {
                $cell.Value = $item[$col]
            }
        }
    }
}
